$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in row 1 (C1: "cuenta" -> "arreglo prestamo",
# D1: "cuenta deposito" -> "cuenta debito"). The old "cuenta" / "cuenta
# deposito" strings are no longer referenced anywhere once this runs, so
# they fall out of the shared-strings table entirely.
$ws.Range("C1").Value = "arreglo prestamo"
$ws.Range("D1").Value = "cuenta debito"

# Move/update the active selection to E1
$ws.Range("E1").Select()
